$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.399.85"
$ws.Range("E2").Value = "  -4.93%  "
$ws.Range("D3").Value = "3.775.64"
$ws.Range("E3").Value = "  -5.32%  "
$ws.Range("D4").Value = "'0.995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").Value = "'578.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.79%  "
$ws.Range("D6").Value = "'162.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").Value = "'0.650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.95%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.727"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").Value = "'0.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").Value = "'51.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.85%  "
$ws.Range("D12").Value = "'0.0000311"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.53%  "
$ws.Range("D13").Value = "'10.92"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "4.337.27"
$ws.Range("E14").Value = "  -5.98%  "
$ws.Range("D15").Value = "3.750.30"
$ws.Range("E15").Value = "  -6.33%  "
$ws.Range("D16").Value = "'20.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'13.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.96%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'1.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.21%  "
$ws.Range("E19").Value = "  -2.84%  "
$ws.Range("D20").Value = "68.846.65"
$ws.Range("E20").Value = "  -5.14%  "
$ws.Range("D21").Value = "'424.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.04%  "
$ws.Range("D22").Value = "'4.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.92%  "
$ws.Range("D23").Value = "'91.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.09%  "
$ws.Range("D24").Value = "'3.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.87%  "
$ws.Range("D25").Value = "'13.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.35%  "
$ws.Range("D26").Value = "'10.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("D27").Value = "'3.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -11.36%  "
$ws.Range("D28").Value = "'5.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "'10.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("D30").Value = "'34.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.35%  "
$ws.Range("D31").Value = "'7.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("D32").Value = "'13.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.41%  "
$ws.Range("D33").Value = "'46.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("D34").Value = "'0.122"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.95%  "
$ws.Range("D35").Value = "'68.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.46%  "
$ws.Range("D36").Value = "0.0₃0955"
$ws.Range("E36").Value = "  +7.75%  "
$ws.Range("D37").Value = "'616.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -8.31%  "
$ws.Range("D38").Value = "'0.415"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.14%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  -0.18%  "
$ws.Range("D41").Value = "'0.141"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D42").Value = "'3.14"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.87%  "
$ws.Range("D43").Value = "'3.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +17.63%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'2.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.78%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0454"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.36%  "
$ws.Range("D46").Value = "'9.60"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -9.72%  "
$ws.Range("D47").Value = "'0.140"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.34%  "
$ws.Range("D48").Value = "'2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -17.55%  "
$ws.Range("D49").Value = "2.795.96"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("D50").Value = "'3.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.13%  "
$ws.Range("D51").Value = "'0.000265"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.99%  "
